# Add files via upload
#
# The uploaded unitType.xlsx catalog grew three new unit-of-measure rows
# (Kilometre, Piece and Thousand pieces), the Hectare row's id switched
# from the zero-padded text "059" to the plain number 59, and the Ton
# row's name was shortened from "Тонна (1000 кг)" to "Тонна". Everything
# else on the sheet is left exactly as it was.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Make room for the new "Километр" row: push the old rows 3-8
# (055.. through 206..) down by one, to rows 4-9.
$ws.Rows.Item(3).Insert()

# New row 3: 8 / Километр / км
$ws.Cells.Item(3, 1).Value = 8
$ws.Cells.Item(3, 2).Value = "Километр"
$ws.Cells.Item(3, 3).Value = "км"

# Row 5 (was "059", Гектар, га) now carries a numeric id instead of the
# zero-padded text id; name/abbreviation are untouched.
$ws.Cells.Item(5, 1).Value = 59

# Row 8 (168, ..., т) - shorten the unit name.
$ws.Cells.Item(8, 2).Value = "Тонна"

# Two brand-new rows appended at the bottom of the table.
$ws.Cells.Item(10, 1).Value = 796
$ws.Cells.Item(10, 2).Value = "Штука"
$ws.Cells.Item(10, 3).Value = "шт"

# Row 11's id ("798") is kept as text, matching how it was originally
# entered in the source catalog (unlike the other new numeric ids above).
$idCell = $ws.Cells.Item(11, 1)
$idCell.NumberFormat = "@"
$idCell.Value = "798"
$idCell.ClearFormats()
$ws.Cells.Item(11, 2).Value = "Тысяча штук"
$ws.Cells.Item(11, 3).Value = "тыс. шт"

# Printer settings recorded for the sheet (portrait, default paper size 9).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection left on B14, as saved in the uploaded file.
$ws.Range("B14").Select()
